# Actualizacion Datos Personales 4 nov
#
# Applies the Nov-4 personal-data update to both roster sheets:
#   - 5ASV (sheet 1): fixes a mistyped phone number and fills in a
#     couple of missing Tutor / Correo_Tutor / Telefono_Tutor entries.
#   - 5AEM (sheet 2): fixes a swapped-letters e-mail alias and fills in
#     a couple of missing Tutor / Correo_Tutor / Telefono_Tutor entries.
#
# Helper: some of the "numbers" in this sheet (phone numbers) are
# stored as *text*, not numeric values (no leading apostrophe visible,
# but the underlying cell type is string). Assigning a numeric-looking
# PowerShell string straight to .Value makes Excel auto-convert it to
# a real number, which would diverge from the workbook's existing
# convention. Instead, drop a text-producing formula ("=""1234""") in
# the target cell, then Copy / PasteSpecial(values-only) it onto
# itself: the paste keeps the String variant type (so the saved cell
# stays t="s") without leaving a formula behind and without touching
# NumberFormat/styles.
function Set-TextValue($ws, $row, $col, [string]$text) {
    $target = $ws.Cells.Item($row, $col)
    $target.Formula = '="' + $text + '"'
    $target.Copy()
    $target.PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = $false
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "5ASV"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("5ASV")

# Row 3 (BONILLA GONZALEZ JOSE EMMANUEL): corrected Tel_Movil + new tutor info
Set-TextValue $ws1 3 6 "2721279750"
$ws1.Cells.Item(3, 8).Value = "MIRIAM GONZÁLEZ AGUILAR"
$ws1.Cells.Item(3, 9).Value = "miriamgonzi84@gmail.com"
Set-TextValue $ws1 3 10 "2721322441"

# Row 20 (MENDOZA GONZALEZ KARLA): new tutor name
$ws1.Cells.Item(20, 8).Value = "AIDA MÉNDOZA GONZÁLEZ"

# Row 34 (TREJO LUENGAS ELIZABETH): new tutor name
$ws1.Cells.Item(34, 8).Value = "JOAQUÍN TREJO MEJÍA"

# ---------------------------------------------------------------
# Sheet "5AEM"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("5AEM")

# Row 18 (PEREZ ROMERO YAIR ANTONIO): corrected e-mail (letters swapped)
$ws2.Cells.Item(18, 5).Value = "yair26prz@gmail.com"

# Row 7 (CARRASCO SANDOVAL CRISTIAN ANTONIO): new tutor info, re-using the
# student's own e-mail/phone as the tutor's contact data
$ws2.Cells.Item(7, 8).Value = "BLANCA ESTELA SANDOVAL DÍAZ"
$ws2.Cells.Item(7, 9).Value = $ws2.Cells.Item(7, 5).Value
Set-TextValue $ws2 7 10 "2721417437"

# Row 37 (VALDERRAMA RODRIGUEZ EMILIO): new tutor info
$ws2.Cells.Item(37, 8).Value = "MARÍA TEREZA RÓDRIGUEZ LOPEZ"
$ws2.Cells.Item(37, 9).Value = "Maytequila133@gmail.com"
Set-TextValue $ws2 37 10 "2721270249"
